$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect date/time-like text columns (B=Date, C=Time) from Excel auto-conversion
$ws.Range("B2:C5").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = 'Saudi Professional League'
$ws.Range("B2").Value = '2025-12-25'
$ws.Range("C2").Value = '14:30:00'
$ws.Range("D2").Value = 'Al Riyadh SC'
$ws.Range("E2").Value = 'Al-Ettifaq'
$ws.Range("F2").Value = 1000
$ws.Range("G2").Value = 1000
$ws.Range("H2").Value = 1.01
$ws.Range("I2").Value = 1.01
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 1000
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 1.02
$ws.Range("O2").Value = 1.09
$ws.Range("P2").Value = 1.01
$ws.Range("Q2").Value = 100
$ws.Range("R2").Value = 1.01
$ws.Range("S2").Value = 200
$ws.Range("T2").Value = 1.99
$ws.Range("U2").Value = 1.01
$ws.Range("V2").Value = 22
$ws.Range("W2").Value = 1.01
$ws.Range("X2").Value = 1000
$ws.Range("Y2").Value = 1.09
$ws.Range("Z2").Value = 1000
$ws.Range("AA2").Value = 1000
$ws.Range("AB2").Value = 1000
$ws.Range("AC2").Value = 1000
$ws.Range("AD2").Value = 1000
$ws.Range("AE2").Value = 1000
$ws.Range("AF2").Value = 1000
$ws.Range("AG2").Value = 1000
$ws.Range("AH2").Value = 1000
$ws.Range("AI2").Value = 1000
$ws.Range("AJ2").Value = 1000
$ws.Range("AK2").Value = 1000
$ws.Range("AL2").Value = 1000
$ws.Range("AM2").Value = 1000
$ws.Range("AN2").Value = 1000
$ws.Range("AO2").Value = 1000

# Row 3
$ws.Range("A3").Value = 'Saudi Professional League'
$ws.Range("B3").Value = '2025-12-25'
$ws.Range("C3").Value = '14:30:00'
$ws.Range("D3").Value = 'NEOM Sports Club'
$ws.Range("E3").Value = 'Al Najma Club'
$ws.Range("F3").Value = 1.05
$ws.Range("G3").Value = 1.06
$ws.Range("H3").Value = 1000
$ws.Range("I3").Value = 1000
$ws.Range("J3").Value = 18.5
$ws.Range("K3").Value = 22
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 1.07
$ws.Range("S3").Value = 1.09
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = 0
$ws.Range("V3").Value = 1.01
$ws.Range("W3").Value = 18
$ws.Range("X3").Value = 1000
$ws.Range("Y3").Value = 1000
$ws.Range("Z3").Value = 1000
$ws.Range("AA3").Value = 1000
$ws.Range("AB3").Value = 1000
$ws.Range("AC3").Value = 1000
$ws.Range("AD3").Value = 1000
$ws.Range("AE3").Value = 1000
$ws.Range("AF3").Value = 1000
$ws.Range("AG3").Value = 1.3
$ws.Range("AH3").Value = 970
$ws.Range("AI3").Value = 1000
$ws.Range("AJ3").Value = 1000
$ws.Range("AK3").Value = 970
$ws.Range("AL3").Value = 1000
$ws.Range("AM3").Value = 1000
$ws.Range("AN3").Value = 1000
$ws.Range("AO3").Value = 1000

# Row 4
$ws.Range("A4").Value = 'Algerian Ligue 1'
$ws.Range("B4").Value = '2025-12-25'
$ws.Range("C4").Value = '15:30:00'
$ws.Range("D4").Value = 'Belouizdad'
$ws.Range("E4").Value = 'ES Setif'
$ws.Range("F4").Value = 1.94
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 8.199999999999999
$ws.Range("I4").Value = 9.4
$ws.Range("J4").Value = 2.56
$ws.Range("K4").Value = 2.7
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 1.36
$ws.Range("N4").Value = 1.59
$ws.Range("O4").Value = 2.62
$ws.Range("P4").Value = 1.14
$ws.Range("Q4").Value = 7.4
$ws.Range("R4").Value = 1.04
$ws.Range("S4").Value = 23
$ws.Range("T4").Value = 5.1
$ws.Range("U4").Value = 1.21
$ws.Range("V4").Value = 1.12
$ws.Range("W4").Value = 2
$ws.Range("X4").Value = 3.85
$ws.Range("Y4").Value = 13
$ws.Range("Z4").Value = 110
$ws.Range("AA4").Value = 1000
$ws.Range("AB4").Value = 3.7
$ws.Range("AC4").Value = 990
$ws.Range("AD4").Value = 90
$ws.Range("AE4").Value = 1000
$ws.Range("AF4").Value = 8.6
$ws.Range("AG4").Value = 28
$ws.Range("AH4").Value = 340
$ws.Range("AI4").Value = 1000
$ws.Range("AJ4").Value = 38
$ws.Range("AK4").Value = 120
$ws.Range("AL4").Value = 1000
$ws.Range("AM4").Value = 440
$ws.Range("AN4").Value = 120
$ws.Range("AO4").Value = 1000

# Row 5
$ws.Range("A5").Value = 'Honduras Liga Nacional'
$ws.Range("B5").Value = '2025-12-25'
$ws.Range("C5").Value = '20:00:00'
$ws.Range("D5").Value = 'Olancho'
$ws.Range("E5").Value = 'Platense FC'
$ws.Range("F5").Value = 2.08
$ws.Range("G5").Value = 2.26
$ws.Range("H5").Value = 3.25
$ws.Range("I5").Value = 3.85
$ws.Range("J5").Value = 3.6
$ws.Range("K5").Value = 4.1
$ws.Range("L5").Value = 1.33
$ws.Range("M5").Value = 1.06
$ws.Range("N5").Value = 4.4
$ws.Range("O5").Value = 1.23
$ws.Range("P5").Value = 2.38
$ws.Range("Q5").Value = 1.67
$ws.Range("R5").Value = 1.55
$ws.Range("S5").Value = 2.6
$ws.Range("T5").Value = 1.55
$ws.Range("U5").Value = 2.42
$ws.Range("V5").Value = 1.36
$ws.Range("W5").Value = 1.8
$ws.Range("X5").Value = 970
$ws.Range("Y5").Value = 1000
$ws.Range("Z5").Value = 1000
$ws.Range("AA5").Value = 1000
$ws.Range("AB5").Value = 1000
$ws.Range("AC5").Value = 1000
$ws.Range("AD5").Value = 1000
$ws.Range("AE5").Value = 1000
$ws.Range("AF5").Value = 1000
$ws.Range("AG5").Value = 1000
$ws.Range("AH5").Value = 1000
$ws.Range("AI5").Value = 1000
$ws.Range("AJ5").Value = 1000
$ws.Range("AK5").Value = 1000
$ws.Range("AL5").Value = 1000
$ws.Range("AM5").Value = 1000
$ws.Range("AN5").Value = 1000
$ws.Range("AO5").Value = 1000

# Remove old row 6 (data shifted up; sheet now ends at row 5)
$ws.Rows("6:6").Delete()